$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H51").Value = 8682.333000000001
$ws.Range("J51").Value = 8500
$ws.Range("L51").Value = 8500
$ws.Range("N51").Value = -9468
$ws.Range("H112").Value = 1323.0385
$ws.Range("J112").Value = 1339.1765
$ws.Range("L112").Value = 4017.5295
$ws.Range("N112").Value = -6233.529500000001
$ws.Range("H129").Value = 2985.7368
$ws.Range("I129").Value = 342.83334
$ws.Range("J129").Value = 4205.5386
$ws.Range("K129").Value = 1028.50002
$ws.Range("L129").Value = 12616.6158
$ws.Range("M129").Value = 3971.49998
$ws.Range("N129").Value = -22616.6158
$ws.Range("H137").Value = 853515
$ws.Range("I137").Value = 2073900.1
$ws.Range("J137").Value = 2943.5757
$ws.Range("K137").Value = 6221700.300000001
$ws.Range("L137").Value = 8830.7271
$ws.Range("M137").Value = -6219150.300000001
$ws.Range("N137").Value = -13930.7271

$ws = $wb.Worksheets("ARM")
$ws.Range("H2").Value = 2482.6667
$ws.Range("I2").Value = 2477.7778
$ws.Range("J2").Value = 2497.3333
$ws.Range("K2").Value = 2477.7778
$ws.Range("L2").Value = 2497.3333
$ws.Range("M2").Value = -2364.7778
$ws.Range("N2").Value = -2723.3333
$ws.Range("H32").Value = 5129.197
$ws.Range("I32").Value = 5929
$ws.Range("K32").Value = 5929
$ws.Range("M32").Value = -5642
$ws.Range("H116").Value = 2482.6667
$ws.Range("I116").Value = 2477.7778
$ws.Range("J116").Value = 2497.3333
$ws.Range("K116").Value = 2477.7778
$ws.Range("L116").Value = 2497.3333
$ws.Range("M116").Value = -183.7777999999998
$ws.Range("N116").Value = -7085.3333
$ws.Range("H132").Value = 2140.851
$ws.Range("I132").Value = 968.88
$ws.Range("J132").Value = 3472.6365
$ws.Range("K132").Value = 2906.64
$ws.Range("L132").Value = 10417.9095
$ws.Range("M132").Value = -376.6399999999999
$ws.Range("N132").Value = -15477.9095

$ws = $wb.Worksheets("BSM")
$ws.Range("H3").Value = 2482.6667
$ws.Range("I3").Value = 2477.7778
$ws.Range("J3").Value = 2497.3333
$ws.Range("K3").Value = 2477.7778
$ws.Range("L3").Value = 2497.3333
$ws.Range("M3").Value = -2363.7778
$ws.Range("N3").Value = -2725.3333

$ws = $wb.Worksheets("CRP")
$ws.Range("H22").Value = 595.3200000000001
$ws.Range("I22").Value = 416.2
$ws.Range("J22").Value = 864
$ws.Range("K22").Value = 416.2
$ws.Range("L22").Value = 864
$ws.Range("M22").Value = -66.19999999999999
$ws.Range("N22").Value = -1564

$ws = $wb.Worksheets("CUL")
$ws.Range("H46").Value = 211.5
$ws.Range("I46").Value = 211.5
$ws.Range("K46").Value = 634.5
$ws.Range("M46").Value = -543.5
$ws.Range("H64").Value = 8001.1665
$ws.Range("I64").Value = 500
$ws.Range("J64").Value = 9501.4
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 28504.2
$ws.Range("M64").Value = -1230
$ws.Range("N64").Value = -29044.2
$ws.Range("H67").Value = 8001.1665
$ws.Range("I67").Value = 500
$ws.Range("J67").Value = 9501.4
$ws.Range("K67").Value = 1500
$ws.Range("L67").Value = 28504.2
$ws.Range("M67").Value = -564
$ws.Range("N67").Value = -30376.2
$ws.Range("H70").Value = 1748
$ws.Range("I70").Value = 1372.6666
$ws.Range("K70").Value = 4117.9998
$ws.Range("M70").Value = -3802.9998
$ws.Range("H73").Value = 1748
$ws.Range("I73").Value = 1372.6666
$ws.Range("K73").Value = 4117.9998
$ws.Range("M73").Value = -3025.9998
$ws.Range("H76").Value = 4500
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 4500
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 13500
$ws.Range("N76").Value = -14266
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 4500
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 4500
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 13500
$ws.Range("N79").Value = -16152
$ws.Range("M79").ClearContents()
$ws.Range("H113").Value = 451.64
$ws.Range("J113").Value = 431.34616
$ws.Range("L113").Value = 1294.03848
$ws.Range("N113").Value = -5634.03848
$ws.Range("H122").Value = 2542.15
$ws.Range("I122").Value = 643.7059
$ws.Range("J122").Value = 3945.348
$ws.Range("K122").Value = 5793.3531
$ws.Range("L122").Value = 35508.132
$ws.Range("M122").Value = -3343.3531
$ws.Range("N122").Value = -40408.132
$ws.Range("H129").Value = 1311.5518
$ws.Range("I129").Value = 993.6875
$ws.Range("J129").Value = 1702.7693
$ws.Range("K129").Value = 2981.0625
$ws.Range("L129").Value = 5108.3079
$ws.Range("M129").Value = 2018.9375
$ws.Range("N129").Value = -15108.3079
$ws.Range("H141").Value = 5925.625
$ws.Range("I141").Value = 5939.231
$ws.Range("J141").Value = 5866.6665
$ws.Range("K141").Value = 17817.693
$ws.Range("L141").Value = 17599.9995
$ws.Range("M141").Value = -12637.693
$ws.Range("N141").Value = -27959.9995

$ws = $wb.Worksheets("GSM")
$ws.Range("H80").Value = 50002320
$ws.Range("I80").Value = 125001000
$ws.Range("J80").Value = 3202
$ws.Range("K80").Value = 125001000
$ws.Range("L80").Value = 3202
$ws.Range("M80").Value = -125000002
$ws.Range("N80").Value = -5198
$ws.Range("H83").Value = 50002320
$ws.Range("I83").Value = 125001000
$ws.Range("J83").Value = 3202
$ws.Range("K83").Value = 625005000
$ws.Range("L83").Value = 16010
$ws.Range("M83").Value = -625000008
$ws.Range("N83").Value = -25994
$ws.Range("H92").Value = 17666.666
$ws.Range("J92").Value = 17666.666
$ws.Range("L92").Value = 17666.666
$ws.Range("N92").Value = -21410.666

$ws = $wb.Worksheets("LTW")
$ws.Range("H16").Value = 861.5714
$ws.Range("I16").Value = 861.5714
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 861.5714
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -691.5714
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets("WVR")
$ws.Range("H6").Value = 2024.4
$ws.Range("J6").Value = 2252.75
$ws.Range("L6").Value = 2252.75
$ws.Range("N6").Value = -2482.75
$ws.Range("H8").Value = 1000003
$ws.Range("I8").Value = 1000003
$ws.Range("K8").Value = 1000003
$ws.Range("M8").Value = -999863
$ws.Range("H9").Value = 10000
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H12").Value = 12069
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 12069
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 12069
$ws.Range("N12").Value = -12353
$ws.Range("M12").ClearContents()
